$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = 'I need to go now.'
$ws.Range("A52").Value = 'I feel good.'
$ws.Range("A53").Value = 'Thank you very much.'
$ws.Range("A54").Value = 'Thanks for your help.'
$ws.Range("A55").Value = 'One like that.'
$ws.Range("A56").Value = 'Be careful.'
$ws.Range("A57").Value = 'I''d like to go home.'
$ws.Range("A58").Value = 'I love you.'
$ws.Range("A59").Value = 'I give up.'
$ws.Range("A60").Value = 'Excuse me.'
$ws.Range("A61").Value = 'Your things are all here.'
$ws.Range("A62").Value = 'This doesn''t work.'
$ws.Range("A63").Value = 'They''ll be right back.'
$ws.Range("A64").Value = 'I''m bored.'
$ws.Range("A65").Value = 'Take it outside.'
$ws.Range("A66").Value = 'Please speak slower.'
$ws.Range("A67").Value = 'Give me the pen.'
$ws.Range("A68").Value = 'If you need my help, please let me know.'
$ws.Range("A69").Value = 'I''m happy.'
$ws.Range("A70").Value = 'I''m going there next year.'
$ws.Range("A71").Value = 'I''m cold.'
$ws.Range("A72").Value = 'Thank you miss.'
$ws.Range("A73").Value = 'Everything is ready.'
$ws.Range("A74").Value = 'I have to wash my clothes.'
$ws.Range("A75").Value = 'That looks old.'
$ws.Range("A76").Value = 'Good afternoon.'
$ws.Range("A77").Value = 'Hurry!'
$ws.Range("A78").Value = 'Please fill out this form.'
$ws.Range("A79").Value = 'Take your time.'
$ws.Range("A80").Value = 'I don''t mind.'
$ws.Range("A81").Value = 'That''s a good school.'
$ws.Range("A82").Value = 'I have one in my car.'
$ws.Range("A83").Value = 'I''ll come back later.'
$ws.Range("A84").Value = 'Tell him that I need to talk to him.'
$ws.Range("A85").Value = 'I''ve never seen that before.'
$ws.Range("A86").Value = 'You''re very smart.'
$ws.Range("A87").Value = 'These books are ours.'
$ws.Range("A88").Value = 'Come here.'
$ws.Range("A89").Value = 'I''ve already seen it.'
$ws.Range("A90").Value = 'I''m just looking.'
$ws.Range("A91").Value = 'Thanks for everything.'
$ws.Range("A92").Value = 'I''d like to use the internet'
$ws.Range("A93").Value = 'Here it is.'
$ws.Range("A94").Value = 'Good idea.'
$ws.Range("A95").Value = 'Please come in.'
$ws.Range("A96").Value = 'I still have a lot to do.'
$ws.Range("A97").Value = 'Happy Birthday.'
$ws.Range("A98").Value = 'I''m getting ready to go out.'
$ws.Range("A99").Value = 'I''ll teach you.'
$ws.Range("A100").Value = 'I feel good.'

$ws.Range("A100").Select()
